# Word-list update: replace 5 "natural" words (col A/B) and 5 "artificial" words
# (col C/D) with new words, then re-sort the table by the "natural" column (A)
# instead of the previous sort-by-"artificial" (col C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column A / B (natural word / natural_freq) replacements, by original row ---
$ws.Cells.Item(10, 1).Value = "צמחים"
$ws.Cells.Item(10, 2).Value = 38

$ws.Cells.Item(11, 1).Value = "עצמות"
$ws.Cells.Item(11, 2).Value = 65

$ws.Cells.Item(17, 1).Value = "ביצים"
$ws.Cells.Item(17, 2).Value = 52

$ws.Cells.Item(19, 1).Value = "פירות"
$ws.Cells.Item(19, 2).Value = 113

$ws.Cells.Item(24, 1).Value = "סוסים"
$ws.Cells.Item(24, 2).Value = 35

# --- Column C / D (artificial word / artificail_freq) replacements, by original row ---
$ws.Cells.Item(6, 3).Value = "פצצות"
$ws.Cells.Item(6, 4).Value = 38

$ws.Cells.Item(11, 3).Value = "טנקים"
$ws.Cells.Item(11, 4).Value = 35

$ws.Cells.Item(16, 3).Value = "פסלים"
$ws.Cells.Item(16, 4).Value = 39

$ws.Cells.Item(20, 3).Value = "ספינה"
$ws.Cells.Item(20, 4).Value = 29

$ws.Cells.Item(21, 3).Value = "רובים"
$ws.Cells.Item(21, 4).Value = 18

# --- Re-sort Table1 by column A (natural) instead of column C (artificial) ---
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("A2:A31"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# --- Update the active selection to match the saved view state ---
$ws.Range("F16").Select()
